# Developed graph for display measurement stamps in time mode.
#
# 1) Typography sheet: the "displayLabel" typography row gets Widget
#    Wildcard Characters set so the numeric graph axis labels
#    ("-., 0123456789") render correctly.
# 2) Translation sheet: add the new text rows (TEXT ID / TYPOGRAPHY NAME /
#    ALIGNMENT / DIRECTION / GB text) used by the new measurement-stamps
#    graph: the graph description, the axis/time-interval labels (7 of
#    them), the graph title, and the generic "<>" label placeholders.

$wb = $excel.ActiveWorkbook
$wsTypo  = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# --- Typography sheet: widen the wildcard characters for displayLabel ---
$wsTypo.Range("H10").Value = "-., 0123456789"

# --- Translation sheet: new rows 226-236 ---
$rows = @(
    @{ Row = 226; B = "SingleUseId286"; C = "displayMeas";  D = "Left";   E = "LTR"; F = "Graph with the measurement values from stamps mode" },
    @{ Row = 227; B = "SingleUseId287"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 1" },
    @{ Row = 228; B = "SingleUseId288"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 2" },
    @{ Row = 229; B = "SingleUseId289"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 3" },
    @{ Row = 230; B = "SingleUseId290"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 4" },
    @{ Row = 231; B = "SingleUseId291"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 5" },
    @{ Row = 232; B = "SingleUseId292"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 6" },
    @{ Row = 233; B = "SingleUseId293"; C = "Default";      D = "Center"; E = "LTR"; F = "Time interval 7" },
    @{ Row = 234; B = "SingleUseId294"; C = "displayMeas";  D = "Left";   E = "LTR"; F = "Graph - Time interval  <value>" },
    @{ Row = 235; B = "SingleUseId295"; C = "displayLabel"; D = "Left";   E = "LTR"; F = "<>" },
    @{ Row = 236; B = "SingleUseId296"; C = "displayLabel"; D = "Right";  E = "LTR"; F = "<>" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $wsTrans.Range("B$n").Value = $r.B
    $wsTrans.Range("C$n").Value = $r.C
    $wsTrans.Range("D$n").Value = $r.D
    $wsTrans.Range("E$n").Value = $r.E
    $wsTrans.Range("F$n").Value = $r.F
}
